# [ADDITIONAL SCRAPING] add a "Player Info" sheet (as the new first sheet)
# with the player's ID/name/batting-hand/bowling-style, and simplify the
# "MATCH_CARD_LINK" column on the "ODI Batting" / "ODI Bowling" sheets into
# a plain "MATCH_CODE" column (just the numeric code instead of the full
# howstat.com scorecard URL).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet in front of the existing sheets.
#    NOTE: sheet references returned by Worksheets.Item(...) track the
#    *position* in the tab strip, so grab the "ODI Batting" / "ODI Bowling"
#    handles AFTER the insert has shifted everything over by one slot.
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Headers
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header look used on the other sheets.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

# Data row for player 4738 - keep the ID as text (matches the source data).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4738"
$playerInfo.Range("B2").Value = "Abu Jayed Chowdhury"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

$playerInfo.Range("A1").Select()

# ------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (store just the code).
# ------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2:D3").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4293"
$battingSheet.Range("D3").Value = "4295"

# ------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (store just the code).
# ------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2:B3").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4293"
$bowlingSheet.Range("B3").Value = "4295"

Write-Host "Sheets: $($wb.Worksheets | ForEach-Object { $_.Name })"
